$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = -10.491
    18 = -11.819
    20 = -12.173
    27 = -13.411
    35 = -12.173
    69 = -10.584
    76 = -12.969
    78 = -12.289
    82 = -11.784
    83 = -13.105
    93 = -10.196
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
